# SimPancake 3000 - add a new todo list
# Updates the "Todo (Game Main)" table with 5 new tasks, fixes two
# existing statuses, logs the new weekly progress numbers on the
# Summary sheet, and refreshes the sheet selections to match the
# state the workbook was left in.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("Summary")
$main    = $wb.Worksheets.Item("Todo (Game Main)")

# ---------------------------------------------------------------
# Todo (Game Main) - "Main" table updates
# ---------------------------------------------------------------

# Two tasks that were finished since the last update
$main.Range("D6").Value  = "C"
$main.Range("D18").Value = "C"

# Five brand-new tasks appended to the table
$main.Range("B20").Value = "Make batter trails look better"
$main.Range("C20").Value = "ie. Mix the trails color up a lil"
$main.Range("D20").Value = "-"

$main.Range("B21").Value = "Improve batter spread"
$main.Range("D21").Value = "C"

$main.Range("B22").Value = "Add dynamic pancake spreed"
$main.Range("D22").Value = "C"

$main.Range("B23").Value = "Make pancakes cook"
$main.Range("D23").Value = "WIP"

$main.Range("B24").Value = "Add AI"

# ---------------------------------------------------------------
# Summary - Progress Log table: freeze last week's live numbers
# and log this week's totals
# ---------------------------------------------------------------

# Week of 19/05 (row 29) - convert the live formulas into a fixed
# historical snapshot
$summary.Range("C29").Value = 24
$summary.Range("D29").Value = 11
$summary.Range("E29").Value = 13

# Week of 26/05 (row 30) - new real counts instead of placeholder 0s
$summary.Range("C30").Value = 28
$summary.Range("D30").Value = 12

# Keep the "Actual % in Week" column chained off the previous row
# now that row 29 is a real data point
$summary.Range("O29").Formula = "=N29-O28"
$summary.Range("O30").Formula = "=N30-O29"
$summary.Range("O31").Formula = "=N31-O30"
$summary.Range("O32").Formula = "=N32-O31"

# Manually-tracked peak % bumped to match the new high
$summary.Range("C11").Value = 0.5

# ---------------------------------------------------------------
# Window / selection state
# ---------------------------------------------------------------

$summary.Activate()
$summary.Application.ActiveWindow.ScrollRow = 4
$summary.Range("B12:D12").Select()

$main.Activate()
$main.Range("C29").Select()

$summary.Activate()
